$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 60, pushing existing rows 60-143 down to 61-144.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new data record.
$ws.Cells.Item(60, 1).Value = 10
$ws.Cells.Item(60, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(60, 3).Value = "La Araucanía"
$ws.Cells.Item(60, 4).Value = 44763
$ws.Cells.Item(60, 5).Value = 9
$ws.Cells.Item(60, 6).Value = 100112031
$ws.Cells.Item(60, 7).Value = "Poroto verde"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 110
$ws.Cells.Item(60, 11).Value = 35000
$ws.Cells.Item(60, 12).Value = 35000
$ws.Cells.Item(60, 13).Value = 35000
$ws.Cells.Item(60, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(60, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(60, 16).Value = 1400
$ws.Cells.Item(60, 17).Value = 25
$ws.Cells.Item(60, 18).Value = "Hortaliza"
